# Scheduled-runner style refresh of market/profit figures across the
# ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets (columns H..N: average prices,
# leve prices and computed profits). Values below are plain data (no
# formulas back these cells), so each change is a direct cell write.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 5353
$ws.Range("I17").Value = 570
$ws.Range("J17").Value = 5584.4355
$ws.Range("K17").Value = 1710
$ws.Range("L17").Value = 16753.3065
$ws.Range("M17").Value = -1542
$ws.Range("N17").Value = -17089.3065
$ws.Range("H43").Value = 1007.5
$ws.Range("I43").Value = 993.3333
$ws.Range("J43").Value = 1050
$ws.Range("K43").Value = 993.3333
$ws.Range("L43").Value = 1050
$ws.Range("M43").Value = -924.3333
$ws.Range("N43").Value = -1188
$ws.Range("H80").Value = 877.6429000000001
$ws.Range("I80").Value = 467.57144
$ws.Range("K80").Value = 1402.71432
$ws.Range("M80").Value = -404.71432
$ws.Range("H83").Value = 877.6429000000001
$ws.Range("I83").Value = 467.57144
$ws.Range("K83").Value = 4208.14296
$ws.Range("M83").Value = 783.8570399999999
$ws.Range("H100").Value = 500
$ws.Range("I100").Value = 500
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 500
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = 41
$ws.Range("N100").ClearContents()
$ws.Range("H132").Value = 34485616
$ws.Range("I132").Value = 37039884
$ws.Range("K132").Value = 111119652
$ws.Range("M132").Value = -111117122
$ws.Range("H135").Value = 1369.5714
$ws.Range("I135").Value = 1084.5
$ws.Range("J135").Value = 1749.6666
$ws.Range("K135").Value = 9760.5
$ws.Range("L135").Value = 15746.9994
$ws.Range("M135").Value = -7225.5
$ws.Range("N135").Value = -20816.9994
$ws.Range("H137").Value = 129289.36
$ws.Range("I137").Value = 179575.8
$ws.Range("J137").Value = 3573.25
$ws.Range("K137").Value = 538727.3999999999
$ws.Range("L137").Value = 10719.75
$ws.Range("M137").Value = -536177.3999999999
$ws.Range("N137").Value = -15819.75
$ws.Range("H138").Value = 2256.5688
$ws.Range("I138").Value = 1009.14636
$ws.Range("J138").Value = 5265.0586
$ws.Range("K138").Value = 3027.43908
$ws.Range("L138").Value = 15795.1758
$ws.Range("M138").Value = 2112.56092
$ws.Range("N138").Value = -26075.1758
$ws.Range("H141").Value = 2486.6667
$ws.Range("I141").Value = 2123.75
$ws.Range("K141").Value = 6371.25
$ws.Range("M141").Value = -1191.25

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 6996810.5
$ws.Range("I45").Value = 11836801
$ws.Range("K45").Value = 11836801
$ws.Range("M45").Value = -11836424
$ws.Range("H63").Value = 1706.2222
$ws.Range("I63").Value = 1415.3334
$ws.Range("K63").Value = 1415.3334
$ws.Range("M63").Value = -729.3334
$ws.Range("H66").Value = 1706.2222
$ws.Range("I66").Value = 1415.3334
$ws.Range("K66").Value = 7076.666999999999
$ws.Range("M66").Value = -3644.666999999999
$ws.Range("H102").Value = 4170014
$ws.Range("J102").Value = 6062.8
$ws.Range("L102").Value = 6062.8
$ws.Range("N102").Value = -9306.799999999999
$ws.Range("H122").Value = 475695.6
$ws.Range("I122").Value = 1864.4849
$ws.Range("J122").Value = 1897188.9
$ws.Range("K122").Value = 5593.4547
$ws.Range("L122").Value = 5691566.699999999
$ws.Range("M122").Value = -3143.4547
$ws.Range("N122").Value = -5696466.699999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2346.182
$ws.Range("I20").Value = 1580.8
$ws.Range("K20").Value = 1580.8
$ws.Range("M20").Value = -1333.8
$ws.Range("H22").Value = 1800.1818
$ws.Range("I22").Value = 1733.3334
$ws.Range("J22").Value = 2101
$ws.Range("K22").Value = 1733.3334
$ws.Range("L22").Value = 2101
$ws.Range("M22").Value = -1560.3334
$ws.Range("N22").Value = -2447
$ws.Range("H50").Value = 40330.668
$ws.Range("J50").Value = 40330.668
$ws.Range("L50").Value = 40330.668
$ws.Range("N50").Value = -41478.668

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 15793.468
$ws.Range("J31").Value = 19089.549
$ws.Range("L31").Value = 19089.549
$ws.Range("N31").Value = -19679.549
$ws.Range("H34").Value = 15793.468
$ws.Range("J34").Value = 19089.549
$ws.Range("L34").Value = 19089.549
$ws.Range("N34").Value = -19493.549
$ws.Range("H122").Value = 4348.7144
$ws.Range("I122").Value = 4165.3335
$ws.Range("K122").Value = 12496.0005
$ws.Range("M122").Value = -10046.0005

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 3352.4
$ws.Range("I69").Value = 1008.5
$ws.Range("J69").Value = 4915
$ws.Range("K69").Value = 3025.5
$ws.Range("L69").Value = 14745
$ws.Range("M69").Value = -2214.5
$ws.Range("N69").Value = -16367
$ws.Range("H72").Value = 3352.4
$ws.Range("I72").Value = 1008.5
$ws.Range("J72").Value = 4915
$ws.Range("K72").Value = 9076.5
$ws.Range("L72").Value = 44235
$ws.Range("M72").Value = -5020.5
$ws.Range("N72").Value = -52347
$ws.Range("H109").Value = 125010690
$ws.Range("I109").Value = 166677260
$ws.Range("K109").Value = 500031780
$ws.Range("M109").Value = -500030740
$ws.Range("H117").Value = 4297.364
$ws.Range("J117").Value = 3903.6
$ws.Range("L117").Value = 11710.8
$ws.Range("N117").Value = -18594.8
$ws.Range("H137").Value = 2725.6
$ws.Range("J137").Value = 4188.4287
$ws.Range("L137").Value = 12565.2861
$ws.Range("N137").Value = -22765.2861

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1629590.8
$ws.Range("I80").Value = 3051080.8
$ws.Range("J80").Value = 5030.857
$ws.Range("K80").Value = 3051080.8
$ws.Range("L80").Value = 5030.857
$ws.Range("M80").Value = -3050082.8
$ws.Range("N80").Value = -7026.857
$ws.Range("H83").Value = 1629590.8
$ws.Range("I83").Value = 3051080.8
$ws.Range("J83").Value = 5030.857
$ws.Range("K83").Value = 15255404
$ws.Range("L83").Value = 25154.285
$ws.Range("M83").Value = -15250412
$ws.Range("N83").Value = -35138.285
$ws.Range("H97").Value = 1833794.1
$ws.Range("I97").Value = 3970588.8
$ws.Range("J97").Value = 2255.8572
$ws.Range("K97").Value = 3970588.8
$ws.Range("L97").Value = 2255.8572
$ws.Range("M97").Value = -3970092.8
$ws.Range("N97").Value = -3247.8572
$ws.Range("H102").Value = 7104387
$ws.Range("I102").Value = 10103440
$ws.Range("J102").Value = 2980689.2
$ws.Range("K102").Value = 10103440
$ws.Range("L102").Value = 2980689.2
$ws.Range("M102").Value = -10101818
$ws.Range("N102").Value = -2983933.2
$ws.Range("H122").Value = 214812.5
$ws.Range("I122").Value = 280287.12
$ws.Range("K122").Value = 840861.36
$ws.Range("M122").Value = -838411.36

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4430.2256
$ws.Range("I40").Value = 3498.5386
$ws.Range("J40").Value = 5103.1113
$ws.Range("K40").Value = 3498.5386
$ws.Range("L40").Value = 5103.1113
$ws.Range("M40").Value = -3362.5386
$ws.Range("N40").Value = -5375.1113
$ws.Range("H68").Value = 2228.2
$ws.Range("I68").Value = 2228.2
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 2228.2
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -1479.2
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 2228.2
$ws.Range("I71").Value = 2228.2
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 11141
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -7397
$ws.Range("N71").ClearContents()
$ws.Range("H122").Value = 5622.8945
$ws.Range("I122").Value = 3732.4546
$ws.Range("K122").Value = 11197.3638
$ws.Range("M122").Value = -8747.363799999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 4850
$ws.Range("H36").Value = 4850
$ws.Range("H40").Value = 19006.25
$ws.Range("I40").Value = 18675
$ws.Range("K40").Value = 18675
$ws.Range("M40").Value = -18526
$ws.Range("H62").Value = 7984.394
$ws.Range("I62").Value = 5001
$ws.Range("J62").Value = 8395.896000000001
$ws.Range("K62").Value = 5001
$ws.Range("L62").Value = 8395.896000000001
$ws.Range("M62").Value = -4377
$ws.Range("N62").Value = -9643.896000000001
$ws.Range("H65").Value = 7984.394
$ws.Range("I65").Value = 5001
$ws.Range("J65").Value = 8395.896000000001
$ws.Range("K65").Value = 25005
$ws.Range("L65").Value = 41979.48
$ws.Range("M65").Value = -21885
$ws.Range("N65").Value = -48219.48
$ws.Range("H122").Value = 2452.818
$ws.Range("I122").Value = 1570.1428
$ws.Range("J122").Value = 3997.5
$ws.Range("K122").Value = 4710.428400000001
$ws.Range("L122").Value = 11992.5
$ws.Range("M122").Value = -2260.428400000001
$ws.Range("N122").Value = -16892.5
$ws.Range("H126").Value = 3472.4443
$ws.Range("I126").Value = 2885.1538
$ws.Range("J126").Value = 4999.4
$ws.Range("K126").Value = 8655.4614
$ws.Range("L126").Value = 14998.2
$ws.Range("M126").Value = -6185.4614
$ws.Range("N126").Value = -19938.2
